# March 24 update 3
# Adds three new columns (renewd / PlanID / iteration) to the sheet, with
# header labels in M1:O1 and constant values for every existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (M1:O1), matching the look/format of the existing headers ---
$ws.Range("L1").Copy() | Out-Null
$ws.Range("M1:O1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"

# --- New data columns for every existing data row (rows 2-23) ---
for ($r = 2; $r -le 23; $r++) {
    $ws.Cells.Item($r, 13).Value = "after"
    $ws.Cells.Item($r, 14).Value = 20141086
    $ws.Cells.Item($r, 15).Value = 2
}
